# Update "想去人数" (want-to-attend count) figures in sheets 展览, 演出, 全部类型
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 633
$ws1.Range("F4").Value = 913
$ws1.Range("F5").Value = 678
$ws1.Range("F6").Value = 812
$ws1.Range("F8").Value = 579
$ws1.Range("F9").Value = 118
$ws1.Range("F11").Value = 601
$ws1.Range("F12").Value = 358
$ws1.Range("F15").Value = 134
$ws1.Range("F16").Value = 320
$ws1.Range("F19").Value = 537
$ws1.Range("F21").Value = 547
$ws1.Range("F23").Value = 573
$ws1.Range("F24").Value = 3

$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F9").Value = 208
$ws2.Range("F11").Value = 15

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F7").Value = 633
$ws4.Range("F8").Value = 913
$ws4.Range("F9").Value = 678
$ws4.Range("F10").Value = 812
$ws4.Range("F12").Value = 579
$ws4.Range("F13").Value = 118
$ws4.Range("F15").Value = 601
$ws4.Range("F18").Value = 358
$ws4.Range("F22").Value = 134
$ws4.Range("F24").Value = 320
$ws4.Range("F27").Value = 208
$ws4.Range("F29").Value = 537
$ws4.Range("F30").Value = 15
$ws4.Range("F34").Value = 547
$ws4.Range("F36").Value = 573
$ws4.Range("F37").Value = 3
